$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D12").Value = -7.203
$ws.Range("D27").Value = -7.843999999999999
$ws.Range("D32").Value = -7.381
$ws.Range("D36").Value = -7.797000000000001
$ws.Range("D38").Value = -7.771999999999998
$ws.Range("D46").Value = -8.156000000000001
$ws.Range("D54").Value = -7.994999999999999
$ws.Range("D55").Value = -8.028
$ws.Range("D56").Value = -8.321
$ws.Range("D67").Value = -7.557
$ws.Range("D69").Value = -7.737
$ws.Range("D72").Value = -7.398000000000001
$ws.Range("D83").Value = -8.038999999999998
$ws.Range("D86").Value = -8.245999999999999
$ws.Range("D91").Value = -7.636
$ws.Range("D93").Value = -7.007000000000001
$ws.Range("D99").Value = -8.074999999999999
$ws.Range("D104").Value = -7.806999999999999
